$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$c = $ws1.Range("C33")
$c.Value = "Rod Weight In Air, Lb"
$c.ReadingOrder = 1
